# Auto-generated: apply cell-value updates from the crypto price refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.416.00'
$ws.Range("D3").Value = '1.912.61'
$ws.Range("E3").Value = '  +1.17%  '
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = '  +0.72%  '
$ws.Range("D5").Value = "'325.37"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("E7").Value = '  +1.44%  '
$ws.Range("D8").Value = "'0.4066"
$ws.Range("E8").Value = '  +0.63%  '
$ws.Range("D9").Value = "'0.08250"
$ws.Range("E9").Value = '  +2.74%  '
$ws.Range("E10").Value = '  +2.32%  '
$ws.Range("D11").Value = "'23.46"
$ws.Range("E11").Value = '  +1.09%  '
$ws.Range("D12").Value = '1.890.95'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = "'6.038"
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("D14").Value = "'7.210"
$ws.Range("E14").Value = '  +2.67%  '
$ws.Range("D15").Value = "'91.14"
$ws.Range("E15").Value = '  +2.38%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = "'1.010"
$ws.Range("E16").Value = '  +0.79%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = "'0.06786"
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").Value = "'17.69"
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").Value = '29.453.62'
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").Value = "'5.624"
$ws.Range("E22").Value = '  +2.09%  '
$ws.Range("E23").Value = '  +1.36%  '
$ws.Range("D24").Value = "'2.194"
$ws.Range("E24").Value = '  +1.53%  '
$ws.Range("D25").Value = '2.117.09'
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("D26").Value = "'6.568"
$ws.Range("E26").Value = '  +10.75%  '
$ws.Range("D27").Value = "'156.69"
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("D28").Value = "'20.05"
$ws.Range("E28").Value = '  +1.75%  '
$ws.Range("E29").Value = '  +1.20%  '
$ws.Range("E30").Value = '  +2.19%  '
$ws.Range("D31").Value = "'1.018"
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("D32").Value = "'0.09551"
$ws.Range("E32").Value = '  +1.46%  '
$ws.Range("D33").Value = "'5.568"
$ws.Range("E33").Value = '  +4.42%  '
$ws.Range("D34").Value = "'3.556"
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("D36").Value = "'0.02283"
$ws.Range("E36").Value = '  +1.74%  '
$ws.Range("E37").Value = '  +1.50%  '
$ws.Range("E38").Value = '  +1.17%  '
$ws.Range("D39").Value = "'8.052"
$ws.Range("E39").Value = '  +2.02%  '
$ws.Range("D40").Value = "'0.5963"
$ws.Range("E40").Value = '  +2.40%  '
$ws.Range("E41").Value = '  +8.00%  '
$ws.Range("D42").Value = "'0.1848"
$ws.Range("E42").Value = '  +1.06%  '
$ws.Range("D43").Value = "'1.279"
$ws.Range("E43").Value = '  -0.58%  '
$ws.Range("D44").Value = "'2.394"
$ws.Range("E44").Value = '  +2.08%  '
$ws.Range("D45").Value = "'0.07614"
$ws.Range("E45").Value = '  -1.04%  '
$ws.Range("D46").Value = "'12.35"
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("D47").Value = "'0.5569"
$ws.Range("E47").Value = '  +1.64%  '
$ws.Range("E48").Value = '  +2.43%  '
$ws.Range("D49").Value = "'117.50"
$ws.Range("E49").Value = '  +4.12%  '
$ws.Range("E50").Value = '  +4.37%  '
$ws.Range("D51").Value = "'72.15"
$ws.Range("E51").Value = '  +1.22%  '
